$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the existing rows 12:13 (dated 45040) down into new rows 14:15,
# preserving all values and formatting (including the date number format).
$ws.Range("A12:T13").Copy() | Out-Null
$ws.Range("A14").PasteSpecial() | Out-Null
$ws.Range("D14").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D15").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Update rows 12:13 with the new (later) report date and revised volume.
$ws.Range("D12").Value = 45049
$ws.Range("D13").Value = 45049
$ws.Range("M13").Value = 60
